$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCaseRun")

# Row 24
$ws.Range("A24").Value = "TB2485446"
$ws.Range("B24").Value = "'11542176"
$ws.Range("C24").Value = "Bond - No Credit"
$ws.Range("D24").Value = "Nil"
$ws.Range("E24").Value = "TC001"

# Row 25
$ws.Range("A25").Value = "TN2485447"
$ws.Range("B25").Value = "'11542189"
$ws.Range("C25").Value = "Personal Auto - Credit"
$ws.Range("D25").Value = "Base"
$ws.Range("E25").Value = "TC002"
